# B1--and-B2-PowerPoint.pptx edit
# 1) Swap the table's style on slide 5 to the new built-in table style.
# 2) Swap the (slide-master) theme colours from the "Integral" / Red Violet
#    palette over to the stock Office theme palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{ABA79C81-C87D-49A7-B2C3-FB569A564271}")
    }
}

# --- 2. Theme colours -------------------------------------------------------
$theme = $p.SlideMaster.Theme
$scheme = $theme.ThemeColorScheme

# Order matches the <a:clrScheme> child order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. Values are VBA-style RGB() integers
# (0x00BBGGRR) for the stock "Office" theme palette.
$officeRgb = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = $officeRgb[$i - 1]
}
